$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update the DataTable values in column A (rows 2-5), keeping the
# "p_palabra" header in A1 untouched.
$ws.Range("A2").Value = "Gallina"
$ws.Range("A3").Value = "Caballo"
$ws.Range("A4").Value = "Teclado"
$ws.Range("A5").Value = "Karting"

# Make "Global" the active sheet/tab and select A4, matching the
# saved view state for execution.
$ws.Activate()
$ws.Range("A4").Select()
